# Update the "Sheet1" worksheet (the active overview sheet) with refreshed
# stats, then move the selection to K11 as in the final saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Update overview values: B2 (Fragments) 2 -> 3, D2 (GDD) 2 -> 3
$ws.Range("B2").Value = 3
$ws.Range("D2").Value = 3

# Move selection to K11 to match final saved cursor position
$ws.Range("K11").Select()
